$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Row, $Values)
    $col = 2  # column B
    foreach ($v in $Values) {
        $ws.Cells.Item($Row, $col).Value = $v
        $col = $col + 1
    }
}

# --- Update table 1 (Min - Max Normalization), rows 3-12, columns B:I ---
Set-RowValues 3  @(3,5,5,2,3,0,2,3)
Set-RowValues 4  @(3,4,5,1,3,1,2,4)
Set-RowValues 5  @(4,3,3,3,2,2,4,2)
Set-RowValues 6  @(4,3,5,2,2,2,2,3)
Set-RowValues 7  @(4,3,5,2,2,2,2,3)
Set-RowValues 8  @(4,4,5,2,2,1,2,3)
Set-RowValues 9  @(4,4,5,2,2,1,2,3)
Set-RowValues 10 @(4,4,5,3,2,1,2,2)
Set-RowValues 11 @(4,4,4,3,2,1,3,2)
Set-RowValues 12 @(4,4,5,2,2,1,2,3)

# Row 13: clear the data (only A13 keeps an empty, unfilled style)
$ws.Range("A13:K13").ClearContents()
$ws.Range("A13").Interior.ColorIndex = 0

# --- Update table 2 (Z-score Normalization), rows 17-26, columns B:I ---
Set-RowValues 17 @(3,5,5,3,3,0,2,2)
Set-RowValues 18 @(3,4,5,3,3,1,2,2)
Set-RowValues 19 @(4,3,5,3,2,2,2,2)
Set-RowValues 20 @(4,3,7,3,2,2,0,2)
Set-RowValues 21 @(4,3,7,3,2,2,0,2)
Set-RowValues 22 @(4,3,6,3,2,2,1,2)
Set-RowValues 23 @(4,3,7,3,2,2,0,2)
Set-RowValues 24 @(4,3,6,3,2,2,1,2)
Set-RowValues 25 @(4,3,6,3,2,2,1,2)
Set-RowValues 26 @(4,3,7,3,2,2,0,2)

# Row 27: clear the data
$ws.Range("A27:K27").ClearContents()
$ws.Range("A27").Interior.ColorIndex = 0

# --- Update the active selection ---
$ws.Range("L26").Select()
